# Update the phone-number sample data on Sayfa1 (A1:A4), center the
# values in column A, and move the active selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: formatted phone number (shared string) -> new sample number
$ws.Range("A1").Value2 = "(531) 111 11 11"

# A2 / A3: raw numeric phone numbers (local + international) -> new sample number
$ws.Range("A2").Value2 = 5311111111
$ws.Range("A3").Value2 = 905311111111

# A4: was a numeric duplicate of A2, now a dash-separated text phone number
$ws.Range("A4").Value2 = "531-111-11-11"

# Center-align the whole column (A1:A4 carry the style in the saved file)
$ws.Columns(1).HorizontalAlignment = -4108   # xlCenter

# Move the selection/active cell to A3
$ws.Range("A3").Select()
